$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the existing "总计" sheet.
#    Clone the "2021-Q4" sheet (same 8-column layout, matching header text
#    and header/index-column styling) so the new sheet keeps the same
#    formatting as its sibling quarter sheets.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(4)
$template   = $wb.Worksheets.Item(3)
$template.Copy($totalSheet)

$q1 = $wb.Worksheets.Item(4)
$q1.Name = "2022-Q1"

# ---------------------------------------------------------------------------
# 2. Overwrite the cloned sheet's data with the 2022-Q1 fund holdings.
#    Columns B-G hold text values in the source data (fund codes keep
#    leading zeros, numeric-looking figures are stored as text) so format
#    those ranges as text before writing.
# ---------------------------------------------------------------------------
$q1.Range("B1:G9").NumberFormat = "@"

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1Data = @(
    @(0, "006102", "浙商丰利增强债券",          "48.75", "47.92", "2.22", "1.0822", 5),
    @(1, "688888", "浙商聚潮产业成长混合",        "8.25",  "93.40", "4.31", "0.3556", 10),
    @(2, "671010", "西部利得策略优选混合A",       "6.79",  "93.07", "5.19", "0.3524", 8),
    @(3, "010381", "浙商智选价值混合A",          "2.92",  "93.43", "4.87", "0.1422", 7),
    @(4, "011060", "西部利得策略优选混合C",       "1.32",  "93.07", "5.19", "0.0685", 8),
    @(5, "007423", "西部利得聚禾灵活配置混合A",    "0.60",  "69.21", "4.68", "0.0281", 4),
    @(6, "007424", "西部利得聚禾灵活配置混合C",    "0.41",  "69.21", "4.68", "0.0192", 4),
    @(7, "010382", "浙商智选价值混合C",          "0.34",  "93.43", "4.87", "0.0166", 7)
)

# Rows 2-5 already exist (copied from the template); rows 6-9 need to be
# created and given the same index-column (A) style as the existing ones.
$q1.Range("A2").Copy()
$q1.Range("A6:A9").PasteSpecial(-4122)  # xlPasteFormats

$row = 2
foreach ($r in $q1Data) {
    $q1.Cells.Item($row, 1).Value = $r[0]
    $q1.Cells.Item($row, 2).Value = $r[1]
    $q1.Cells.Item($row, 3).Value = $r[2]
    $q1.Cells.Item($row, 4).Value = $r[3]
    $q1.Cells.Item($row, 5).Value = $r[4]
    $q1.Cells.Item($row, 6).Value = $r[5]
    $q1.Cells.Item($row, 7).Value = $r[6]
    $q1.Cells.Item($row, 8).Value = $r[7]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3. Update the "总计" (totals) sheet: insert a new row for 2022-Q1 above the
#    existing 2021-Q4 row, pushing the older rows down, and renumber index
#    column A.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(5)
$total.Rows.Item(2).Insert()

# The blank inserted row picks up stray bold formatting on B:D - clear it so
# the new data row matches the plain (unstyled) look of the other data rows.
$total.Range("B2:D2").ClearFormats()

# Column A keeps the bold/bordered index style; copy it from the row that
# was just pushed down into row 3 (must copy AFTER the insert, since Insert
# invalidates any earlier clipboard contents).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 2.06

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
